# The filenames/notes for the "44" and "45" supplementary programs were
# listed in the wrong order (authexp/hindex swapped). Fix the order so
# that 44 = hindex and 45 = authexp.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A19").Value = "./programs/44_supplementary_hindex.R"
$ws.Range("B19").Value = "Re-computes the h-index per author and year"
$ws.Range("A20").Value = "./programs/45_supplementary_authexp.R"
$ws.Range("B20").Value = "This program primarily handles computing academic age of author"

# Align the Note column's cell formatting with the Filename column on the
# rows where it had drifted.
$ws.Range("A2").Copy()
$ws.Range("B2").PasteSpecial(-4122)
$ws.Range("A3").Copy()
$ws.Range("B3").PasteSpecial(-4122)

$ws.Range("A12").Copy()
$ws.Range("B12").PasteSpecial(-4122)
$ws.Range("A13").Copy()
$ws.Range("B13").PasteSpecial(-4122)
$ws.Range("A14").Copy()
$ws.Range("B14").PasteSpecial(-4122)

$ws.Range("A17").Copy()
$ws.Range("B17:B31").PasteSpecial(-4122)

$excel.CutCopyMode = $false

# Move the active selection to where the editor last left it.
[void]$ws.Range("C29").Select()
